$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("T2").Value = 1
$ws.Range("V2").Value = 10
$ws.Range("W2").Value = 5000
$ws.Range("X2").Value = 1500

[void]$ws.Range("S6").Select()
